$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find & Replace).
$d.Content.Find.Execute("2024-12-04 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-12-05 Thursday", 2) | Out-Null

# Update the division problems inside the table.
# Cells are addressed directly by (row, column) rather than via a global
# text Find & Replace, because some new values collide with old values
# elsewhere in the table (e.g. "86÷6=" is both a pre-edit value and a
# post-edit value), which would make a naive global replace unsafe.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "87÷5="
$t.Cell(1, 2).Range.Text  = "62÷7="
$t.Cell(1, 3).Range.Text  = "74÷5="
$t.Cell(1, 4).Range.Text  = "14÷5="
$t.Cell(1, 5).Range.Text  = "35÷7="

$t.Cell(5, 1).Range.Text  = "83÷3="
$t.Cell(5, 2).Range.Text  = "62÷2="
$t.Cell(5, 3).Range.Text  = "37÷8="
$t.Cell(5, 4).Range.Text  = "49÷7="
$t.Cell(5, 5).Range.Text  = "86÷6="

$t.Cell(9, 1).Range.Text  = "16÷4="
$t.Cell(9, 2).Range.Text  = "96÷9="
$t.Cell(9, 3).Range.Text  = "98÷7="
$t.Cell(9, 4).Range.Text  = "33÷7="
$t.Cell(9, 5).Range.Text  = "59÷6="

$t.Cell(13, 1).Range.Text = "10÷8="
$t.Cell(13, 2).Range.Text = "86÷4="
$t.Cell(13, 3).Range.Text = "94÷2="
$t.Cell(13, 4).Range.Text = "26÷7="
$t.Cell(13, 5).Range.Text = "54÷8="

$t.Cell(17, 1).Range.Text = "30÷5="
$t.Cell(17, 2).Range.Text = "50÷4="
$t.Cell(17, 3).Range.Text = "73÷6="
$t.Cell(17, 4).Range.Text = "67÷3="
$t.Cell(17, 5).Range.Text = "57÷8="
